$d = $word.ActiveDocument

# 1) Remove the _GoBack bookmark from its current location (end of the JSX paragraph);
#    it will be re-inserted at the new location below.
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Insert the three new paragraphs (spacer, new question, new answer) after the
#    last paragraph in the document (the "JSX is a syntax extension..." answer).
$insertPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="7"/><w:numPr><w:numId w:val="0"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:left="360" w:leftChars="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="202122"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="7"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="720" w:leftChars="0" w:hanging="360" w:firstLineChars="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr><w:t>Explain how lists work in React.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="7"/><w:numPr><w:numId w:val="0"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:left="360" w:leftChars="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="202122"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="202122"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr><w:t>We create lists in React as we do in regular JavaScript. Lists display data in an ordered format. The traversal of lists is done using the map() function</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$null = $insertPoint.InsertXML($xml)
